$d = $word.ActiveDocument

# wdHeaderFooterPrimary = 1
$footer = $d.Sections(1).Footers(1)

# --- Footer: bump version number "Version 1.8.x" -> "Version 2.0.x" ---
$footer.Range.Find.Execute("ersion 1.8.x", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "ersion 2.0.x", 2)

# --- Footer: update the cached "Last update" date ---
$footer.Range.Find.Execute("2020-10-05", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2024-03-08", 2)
